$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Story Board")

$ws.Range("B6").Value = "Understand how to upload code in iterations on GitHub"
$ws.Range("C6").Value = "Dan"

$ws.Columns.Item(2).ColumnWidth = 43.6640625
$ws.Columns.Item(3).ColumnWidth = 43.6640625

$ws.Range("C6").Select()
